# Updated symbol list on Thu Jan 12 09:22:57 UTC 2023 with GitHub Actions
# Apply refreshed Price (D) and Volume(1h) (E) values for the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "283.41";    "E2"  = "1.99%"
    "D3"  = "28.46";     "E3"  = "4.33%"
    "D4"  = "5.047";     "E4"  = "3.47%"
    "D5"  = "0.06497";   "E5"  = "1.08%"
    "D6"  = "7.217";     "E6"  = "3.39%"
    "D7"  = "1.426";     "E7"  = "19.97%"
    "D8"  = "0.9172";    "E8"  = "3.78%"
    "D9"  = "0.1540";    "E9"  = "-1.79%"
    "D10" = "0.06582";   "E10" = "28.27%"
    "D11" = "0.07603";   "E11" = "2.51%"
    "D12" = "0.02787";   "E12" = "-3.28%"
    "D13" = "0.08966";   "E13" = "-0.13%"
    "D14" = "0.001583";  "E14" = "0.82%"
    "D15" = "0.0006344"; "E15" = "-0.87%"
    "D16" = "0.006071";  "E16" = "-1.15%"
    "D17" = "3.453";     "E17" = "-0.88%"
    "D18" = "3.365";     "E18" = "1.61%"
    "E19" = "-1.42%"
    "E20" = "1.19%"
    "D21" = "0.1341";    "E21" = "-0.63%"
    "D22" = "3.979";     "E22" = "1.96%"
    "E23" = "2.94%"
    "D24" = "0.04442";   "E24" = "0.35%"
    "D25" = "0.001182";  "E25" = "0.41%"
    "D26" = "0.004464";  "E26" = "15.34%"
    "D28" = "0.0001200"; "E28" = "1.78%"
    "E29" = "-15.73%"
    "D40" = "0.04116";   "E40" = "-0.66%"
    "D41" = "0.006677";  "E41" = "-2.46%"
    "D42" = "0.1232";    "E42" = "4.95%"
    "D43" = "0.002050";  "E43" = "5.22%"
    "D44" = "0.01244";   "E44" = "9.12%"
    "D45" = "0.00005401"; "E45" = "1.83%"
    "E46" = "-0.07%"
}

# Values are stored as plain text in the sheet (e.g. "283.41", "1.99%"),
# not as numbers. A leading apostrophe forces Excel to keep the literal
# text instead of auto-converting to a number/percentage.
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
